# "Complete test for table_CAPM macro"
#
# The sheet lists one test per row (Test name | Description | macro),
# roughly alphabetised by column A. The placeholder row for the CAPM
# table test ("table_CAPM1" / "Test CAPM table" / "table_CAPM_test1")
# is removed and replaced by two fully specified rows appended at the
# bottom of the table:
#   table_CAPM1 | Test CAPM table with scale=252, digits=4 | table_CAPM_test1
#   table_CAPM2 | Test CAPM table with scale=1, digits=6   | table_CAPM_test2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old placeholder row (row 53: table_CAPM1 / Test CAPM table / table_CAPM_test1).
$ws.Rows(53).Delete()

# After the deletion the table now spans rows 2-74 (row 1 is the header),
# so the two new rows land at 74 and 75.
$ws.Range("A74").Value = "table_CAPM1"
$ws.Range("B74").Value = "Test CAPM table with scale=252, digits=4"
$ws.Range("C74").Value = "table_CAPM_test1"

$ws.Range("A75").Value = "table_CAPM2"
$ws.Range("B75").Value = "Test CAPM table with scale=1, digits=6"
$ws.Range("C75").Value = "table_CAPM_test2"

# Match the workbook's saved selection state (last touched cell).
$ws.Range("C75").Select()
